# Fruta / hortaliza, semanal
# Insert two new weekly data rows (Mapocho Venta Directa de Santiago - Esparragos)
# right after the header block of existing "Sin especificar" rows at row 20,
# pushing the remaining historical rows down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at row 20; existing rows 20:55 shift down to 22:57
$ws.Rows("20:21").Insert()

# New row 20 - Banquete
$ws.Range("A20").Value = 12
$ws.Range("B20").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C20").Value = "Metropolitana"
$ws.Range("D20").Value = 44495
$ws.Range("E20").Value = 13
$ws.Range("F20").Value = 300000000
$ws.Range("G20").Value = "Espárragos"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Banquete"
$ws.Range("J20").Value = 300
$ws.Range("K20").Value = 1200
$ws.Range("L20").Value = 1200
$ws.Range("M20").Value = 1200
$ws.Range("N20").Value = "$/kilo"
$ws.Range("O20").Value = "Región Metropolitana"
$ws.Range("P20").Value = 1200
$ws.Range("Q20").Value = 1
$ws.Range("R20").Value = "Hortaliza"

# New row 21 - Primera
$ws.Range("A21").Value = 12
$ws.Range("B21").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C21").Value = "Metropolitana"
$ws.Range("D21").Value = 44495
$ws.Range("E21").Value = 13
$ws.Range("F21").Value = 300000000
$ws.Range("G21").Value = "Espárragos"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 350
$ws.Range("K21").Value = 1000
$ws.Range("L21").Value = 1000
$ws.Range("M21").Value = 1000
$ws.Range("N21").Value = "$/kilo"
$ws.Range("O21").Value = "Región Metropolitana"
$ws.Range("P21").Value = 1000
$ws.Range("Q21").Value = 1
$ws.Range("R21").Value = "Hortaliza"
